# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The upstream data-generation script recomputes the "K" column (column G,
# header "K") for every trade row on Sheet1 from the freshly regenerated
# strike-count source data. This updates each row's K value in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (column G) value, as produced by the regenerated
# save_data pipeline.
$kValues = [ordered]@{
    2  = 0
    3  = 2
    4  = 1
    5  = 0
    6  = 2
    7  = 1
    8  = 2
    9  = 1
    10 = 1
    11 = 0
    12 = 0
    13 = 1
    16 = 2
    17 = 0
    18 = 1
    19 = 0
    20 = 1
    21 = 0
    22 = 1
    24 = 1
    25 = 1
    26 = 3
    27 = 0
    29 = 1
    30 = 0
    31 = 1
    32 = 1
    33 = 0
    34 = 2
    35 = 0
    36 = 2
    37 = 2
    38 = 0
    39 = 2
    40 = 1
    41 = 0
    42 = 2
    43 = 0
    44 = 0
    45 = 1
    46 = 3
    47 = 0
    48 = 2
    49 = 1
    50 = 1
    51 = 0
    52 = 0
    53 = 0
    54 = 0
    55 = 1
    56 = 0
    57 = 1
    58 = 1
    59 = 3
    61 = 2
    62 = 1
    63 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
